# Refresh the crypto price/volume snapshot (Price column D, Volume(1h) column E)
# to match the latest scrape. Some rows also swap position with their neighbour
# (ShibaInu <-> WrappedBTC at rows 16/17, Mantle <-> RenderToken at rows 47/48),
# so Coin (B) and Link (C) are updated there too.
#
# Several new Price values are single numeric tokens (e.g. "0.9995", "1.0000")
# that Excel would otherwise auto-convert to a number, collapsing formatting
# such as trailing zeros. Force those specific cells to Text ("@") first so
# the literal string is preserved, matching the source data which stores all
# Price values as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.369.42"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "1.877.48"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7135"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.03"
$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3116"
$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07713"
$ws.Range("E9").Value = "  -2.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.14"
$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08372"
$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("D12").Value = "1.909.94"
$ws.Range("E12").Value = "  +2.04%  "

$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("E14").Value = "  -1.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.79"

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "29.371.55"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008324"
$ws.Range("E17").Value = "  +6.26%  "

$ws.Range("E18").Value = "  +1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.55"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "2.133.50"
$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.22"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.885"
$ws.Range("E23").Value = "  -1.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.0000"

$ws.Range("E25").Value = "  +1.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.79"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.026"
$ws.Range("E27").Value = "  +0.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.56"
$ws.Range("E28").Value = "  +1.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.505"
$ws.Range("E29").Value = "  +0.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.406"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.338"
$ws.Range("E31").Value = "  +5.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.291"
$ws.Range("E32").Value = "  -4.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05255"
$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.922"
$ws.Range("E34").Value = "  -0.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7573"
$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.175"
$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01864"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.717"
$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("D40").Value = "1.162.48"
$ws.Range("E40").Value = "  -0.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.350"
$ws.Range("E41").Value = "  +3.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.67"
$ws.Range("E42").Value = "  +1.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8893"
$ws.Range("E43").Value = "  -1.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.78"
$ws.Range("E44").Value = "  +2.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "2.031.94"
$ws.Range("E46").Value = "  +0.82%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.799"
$ws.Range("E47").Value = "  +0.80%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5199"
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.410"
$ws.Range("E49").Value = "  +1.41%  "

$ws.Range("E50").Value = "  -0.25%  "

$ws.Range("E51").Value = "  +0.54%  "
